$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 71, shifting all existing rows (71..162) down to (72..163).
$ws.Rows("71:71").Insert()

# Populate the newly inserted row 71 with the new weekly data entry.
$ws.Cells.Item(71, 1).Value = 3
$ws.Cells.Item(71, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(71, 3).Value = "Coquimbo"
$ws.Cells.Item(71, 4).Value = 44638
$ws.Cells.Item(71, 5).Value = 5
$ws.Cells.Item(71, 6).Value = 100112052
$ws.Cells.Item(71, 7).Value = "Albahaca"
$ws.Cells.Item(71, 8).Value = "Sin especificar"
$ws.Cells.Item(71, 9).Value = "Primera"
$ws.Cells.Item(71, 10).Value = 70
$ws.Cells.Item(71, 11).Value = 4000
$ws.Cells.Item(71, 12).Value = 4000
$ws.Cells.Item(71, 13).Value = 4000
$ws.Cells.Item(71, 14).Value = "`$/docena de matas"
$ws.Cells.Item(71, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(71, 16).Value = 667
$ws.Cells.Item(71, 17).Value = 6
$ws.Cells.Item(71, 18).Value = "Hortaliza"
